$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new "Erreichte Punkte" column (C) values
$ws.Range("C2").Value = 9
$ws.Range("C4").Value = 5
$ws.Range("C5").Value = 5
$ws.Range("C7").Value = 5
$ws.Range("C8").Value = 10
$ws.Range("C10").Value = 5

# Slightly adjusted column widths (as stored in the worksheet XML, Excel's
# ColumnWidth property is offset by 5/6 of a character from the stored value)
$ws.Columns.Item(2).ColumnWidth = 14.33203125 - (5/6)
$ws.Columns.Item(3).ColumnWidth = 23.109375 - (5/6)
$ws.Columns.Item(4).ColumnWidth = 17.44140625 - (5/6)

# Update the active selection to C8, matching the saved view state
$ws.Range("C8").Select()
